$wb = $excel.ActiveWorkbook

# --- 1. Update the "Ready for handoff" status text to "In Translation" ---
# This shared string is referenced by the status cells on all three sheets;
# setting each cell's value updates every occurrence.
$overview = $wb.Worksheets.Item(1)
$overview.Range("E2").Value2 = "In Translation"
$overview.Range("F2").Value2 = "In Translation"
$overview.Range("E3").Value2 = "In Translation"
$overview.Range("F3").Value2 = "In Translation"

$zhcn = $wb.Worksheets.Item(2)
$zhcn.Range("C2").Value2 = "In Translation"
$zhcn.Range("C3").Value2 = "In Translation"

$dede = $wb.Worksheets.Item(3)
$dede.Range("C2").Value2 = "In Translation"
$dede.Range("C3").Value2 = "In Translation"

# --- 2. Narrow the status columns (zh-cn / de-de) on all three sheets ---
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
